# The sentence in the document reads:
#   Step 16 The system displays the “Exit” or “Close another admission” prompt.
# The edit removes the "Exit" or  portion (run-splitting aside, the visible
# text drops "Exit” or “" and leaves only the opening curly quote before
# "Close another admission").
#
# Result:
#   Step 16 The system displays the “Close another admission” prompt.

$d = $word.ActiveDocument

$find = "Exit" + [char]0x201D + " or " + [char]0x201C
$replace = ""

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2) | Out-Null
